$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2233670
$ws.Range("I17").Value = 1213.3334
$ws.Range("J17").Value = 3628955.5
$ws.Range("K17").Value = 3640.0002
$ws.Range("L17").Value = 10886866.5
$ws.Range("M17").Value = -3472.0002
$ws.Range("N17").Value = -10887202.5

$ws.Range("H103").Value = 699.8570999999999
$ws.Range("I103").Value = 800
$ws.Range("K103").Value = 2400
$ws.Range("M103").Value = -1814

$ws.Range("H106").Value = 18098.04
$ws.Range("I106").Value = 3497.8667
$ws.Range("K106").Value = 3497.8667
$ws.Range("M106").Value = -2866.8667

$ws.Range("H115").Value = 229.8
$ws.Range("I115").Value = 229.8
$ws.Range("K115").Value = 689.4000000000001
$ws.Range("M115").Value = 877.5999999999999

$ws.Range("H118").Value = 777.63635
$ws.Range("I118").Value = 776.25
$ws.Range("J118").Value = 781.3333
$ws.Range("K118").Value = 2328.75
$ws.Range("L118").Value = 2343.9999
$ws.Range("M118").Value = -671.75
$ws.Range("N118").Value = -5657.9999

$ws.Range("H132").Value = 1814.0286
$ws.Range("I132").Value = 1864.2142
$ws.Range("J132").Value = 1613.2858
$ws.Range("K132").Value = 5592.642599999999
$ws.Range("L132").Value = 4839.857400000001
$ws.Range("M132").Value = -3062.642599999999
$ws.Range("N132").Value = -9899.857400000001

$ws.Range("H136").Value = 43890
$ws.Range("J136").Value = 43890
$ws.Range("L136").Value = 43890
$ws.Range("N136").Value = -54090

$ws.Range("H137").Value = 14507.182
$ws.Range("I137").Value = 5698.737
$ws.Range("J137").Value = 26461.5
$ws.Range("K137").Value = 17096.211
$ws.Range("L137").Value = 79384.5
$ws.Range("M137").Value = -14546.211
$ws.Range("N137").Value = -84484.5

$ws.Range("H138").Value = 2840.1177
$ws.Range("I138").Value = 761.2083
$ws.Range("J138").Value = 4688.037
$ws.Range("K138").Value = 2283.6249
$ws.Range("L138").Value = 14064.111
$ws.Range("M138").Value = 2856.3751
$ws.Range("N138").Value = -24344.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2019.3062
$ws.Range("I2").Value = 1349.1177
$ws.Range("J2").Value = 3538.4
$ws.Range("K2").Value = 1349.1177
$ws.Range("L2").Value = 3538.4
$ws.Range("M2").Value = -1236.1177
$ws.Range("N2").Value = -3764.4

$ws.Range("H45").Value = 11868.363
$ws.Range("I45").Value = 13815.444
$ws.Range("J45").Value = 3106.5
$ws.Range("K45").Value = 13815.444
$ws.Range("L45").Value = 3106.5
$ws.Range("M45").Value = -13438.444
$ws.Range("N45").Value = -3860.5

$ws.Range("H61").Value = 5314.6665
$ws.Range("I61").Value = 3157.5715
$ws.Range("K61").Value = 3157.5715
$ws.Range("M61").Value = -2945.5715

$ws.Range("H88").Value = 2359.8
$ws.Range("I88").Value = 2199.5
$ws.Range("J88").Value = 2466.6667
$ws.Range("K88").Value = 2199.5
$ws.Range("L88").Value = 2466.6667
$ws.Range("M88").Value = -1793.5
$ws.Range("N88").Value = -3278.6667

$ws.Range("H91").Value = 2359.8
$ws.Range("I91").Value = 2199.5
$ws.Range("J91").Value = 2466.6667
$ws.Range("K91").Value = 2199.5
$ws.Range("L91").Value = 2466.6667
$ws.Range("M91").Value = -795.5
$ws.Range("N91").Value = -5274.6667

$ws.Range("H116").Value = 2019.3062
$ws.Range("I116").Value = 1349.1177
$ws.Range("J116").Value = 3538.4
$ws.Range("K116").Value = 1349.1177
$ws.Range("L116").Value = 3538.4
$ws.Range("M116").Value = 944.8823
$ws.Range("N116").Value = -8126.4

$ws.Range("H132").Value = 2888.0896
$ws.Range("I132").Value = 2433.509
$ws.Range("K132").Value = 7300.527
$ws.Range("M132").Value = -4770.527

$ws.Range("H136").Value = 5314.6665
$ws.Range("I136").Value = 3157.5715
$ws.Range("K136").Value = 9472.7145
$ws.Range("M136").Value = -6922.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2019.3062
$ws.Range("I3").Value = 1349.1177
$ws.Range("J3").Value = 3538.4
$ws.Range("K3").Value = 1349.1177
$ws.Range("L3").Value = 3538.4
$ws.Range("M3").Value = -1235.1177
$ws.Range("N3").Value = -3766.4

$ws.Range("H22").Value = 85
$ws.Range("I22").Value = 85.625
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 85.625
$ws.Range("L22").Value = 80
$ws.Range("M22").Value = 87.375
$ws.Range("N22").Value = -426

$ws.Range("H103").Value = 35905.75
$ws.Range("J103").Value = 35905.75
$ws.Range("L103").Value = 35905.75
$ws.Range("N103").Value = -38249.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4076.121
$ws.Range("I31").Value = 3519.2942
$ws.Range("J31").Value = 4667.75
$ws.Range("K31").Value = 3519.2942
$ws.Range("L31").Value = 4667.75
$ws.Range("M31").Value = -3224.2942
$ws.Range("N31").Value = -5257.75

$ws.Range("H34").Value = 4076.121
$ws.Range("I34").Value = 3519.2942
$ws.Range("J34").Value = 4667.75
$ws.Range("K34").Value = 3519.2942
$ws.Range("L34").Value = 4667.75
$ws.Range("M34").Value = -3317.2942
$ws.Range("N34").Value = -5071.75

$ws.Range("H58").Value = 3911.5908
$ws.Range("I58").Value = 2230.111
$ws.Range("J58").Value = 11478.25
$ws.Range("K58").Value = 2230.111
$ws.Range("L58").Value = 11478.25
$ws.Range("M58").Value = -2027.111
$ws.Range("N58").Value = -11884.25

$ws.Range("H99").Value = 8110.32
$ws.Range("I99").Value = 3661.8635
$ws.Range("J99").Value = 11605.536
$ws.Range("K99").Value = 3661.8635
$ws.Range("L99").Value = 11605.536
$ws.Range("M99").Value = -2163.8635
$ws.Range("N99").Value = -14601.536

$ws.Range("H105").Value = 2055.5
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

$ws.Range("H122").Value = 1637.25
$ws.Range("I122").Value = 1756.4286
$ws.Range("J122").Value = 1470.4
$ws.Range("K122").Value = 5269.2858
$ws.Range("L122").Value = 4411.200000000001
$ws.Range("M122").Value = -2819.2858
$ws.Range("N122").Value = -9311.200000000001

$ws.Range("H126").Value = 8110.32
$ws.Range("I126").Value = 3661.8635
$ws.Range("J126").Value = 11605.536
$ws.Range("K126").Value = 10985.5905
$ws.Range("L126").Value = 34816.608
$ws.Range("M126").Value = -8515.5905
$ws.Range("N126").Value = -39756.608

$ws.Range("H132").Value = 25134.016
$ws.Range("I132").Value = 16308.256
$ws.Range("K132").Value = 48924.768
$ws.Range("M132").Value = -46394.768

$ws.Range("H134").Value = 2281.717
$ws.Range("I134").Value = 1820.0426
$ws.Range("J134").Value = 5898.1665
$ws.Range("K134").Value = 5460.1278
$ws.Range("L134").Value = 17694.4995
$ws.Range("M134").Value = -2925.1278
$ws.Range("N134").Value = -22764.4995

$ws.Range("H136").Value = 3911.5908
$ws.Range("I136").Value = 2230.111
$ws.Range("J136").Value = 11478.25
$ws.Range("K136").Value = 6690.333
$ws.Range("L136").Value = 34434.75
$ws.Range("M136").Value = -4140.333
$ws.Range("N136").Value = -39534.75

$ws.Range("H141").Value = 542548.8
$ws.Range("I141").Value = 312647.5
$ws.Range("J141").Value = 657499.5
$ws.Range("K141").Value = 312647.5
$ws.Range("L141").Value = 657499.5
$ws.Range("M141").Value = -307467.5
$ws.Range("N141").Value = -667859.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4894
$ws.Range("J32").Value = 4894
$ws.Range("L32").Value = 14682
$ws.Range("N32").Value = -15248

$ws.Range("H39").Value = 5770.2144
$ws.Range("J39").Value = 8103.6665
$ws.Range("L39").Value = 24310.9995
$ws.Range("N39").Value = -24898.9995

$ws.Range("H55").Value = 4434.5
$ws.Range("I55").Value = 2473.25
$ws.Range("J55").Value = 5742
$ws.Range("K55").Value = 7419.75
$ws.Range("L55").Value = 17226
$ws.Range("M55").Value = -7242.75
$ws.Range("N55").Value = -17580

$ws.Range("H131").Value = 9461.529
$ws.Range("I131").Value = 1337.5
$ws.Range("J131").Value = 11961.23
$ws.Range("K131").Value = 4012.5
$ws.Range("L131").Value = 35883.69
$ws.Range("M131").Value = 1027.5
$ws.Range("N131").Value = -45963.69

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 12717
$ws.Range("I132").Value = 9570.786
$ws.Range("J132").Value = 27399.334
$ws.Range("K132").Value = 28712.358
$ws.Range("L132").Value = 82198.00199999999
$ws.Range("M132").Value = -26182.358
$ws.Range("N132").Value = -87258.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 818.6667
$ws.Range("I22").Value = 782
$ws.Range("J22").Value = 892
$ws.Range("K22").Value = 782
$ws.Range("L22").Value = 892
$ws.Range("M22").Value = -487
$ws.Range("N22").Value = -1482

$ws.Range("H27").Value = 818.6667
$ws.Range("I27").Value = 782
$ws.Range("J27").Value = 892
$ws.Range("K27").Value = 782
$ws.Range("L27").Value = 892
$ws.Range("M27").Value = -675
$ws.Range("N27").Value = -1106

$ws.Range("H122").Value = 4006.5833
$ws.Range("I122").Value = 4009
$ws.Range("J122").Value = 3994.5
$ws.Range("K122").Value = 12027
$ws.Range("L122").Value = 11983.5
$ws.Range("M122").Value = -9577
$ws.Range("N122").Value = -16883.5

$ws.Range("H132").Value = 4625.2188
$ws.Range("I132").Value = 3777.3333
$ws.Range("J132").Value = 5373.353
$ws.Range("K132").Value = 11331.9999
$ws.Range("L132").Value = 16120.059
$ws.Range("M132").Value = -8801.999899999999
$ws.Range("N132").Value = -21180.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14062.5
$ws.Range("I45").Value = 11812
$ws.Range("J45").Value = 16313
$ws.Range("K45").Value = 11812
$ws.Range("L45").Value = 16313
$ws.Range("M45").Value = -11321
$ws.Range("N45").Value = -17295

$ws.Range("H107").Value = 2527879
$ws.Range("I107").Value = 2246.0833
$ws.Range("J107").Value = 5558638.5
$ws.Range("K107").Value = 6738.249899999999
$ws.Range("L107").Value = 16675915.5
$ws.Range("M107").Value = -4818.249899999999
$ws.Range("N107").Value = -16679755.5

$ws.Range("H122").Value = 5001.1387
$ws.Range("I122").Value = 3757.6667
$ws.Range("J122").Value = 8731.556
$ws.Range("K122").Value = 11273.0001
$ws.Range("L122").Value = 26194.668
$ws.Range("M122").Value = -8823.000100000001
$ws.Range("N122").Value = -31094.668

$ws.Range("H132").Value = 16124.957
$ws.Range("I132").Value = 9775.395500000001
$ws.Range("J132").Value = 29978.545
$ws.Range("K132").Value = 29326.1865
$ws.Range("L132").Value = 89935.63499999999
$ws.Range("M132").Value = -26796.1865
$ws.Range("N132").Value = -94995.63499999999
